$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "SamplesTab" row's query (B3) so the Tumor column resolves
# the sample's tumor status directly instead of via the collected list,
# matching the updated Phs accession / Study / Gender filter test cases.
$newTumorQuery = @'

MATCH (s:study)<--(p:participant)<--(samp:sample)
WHERE s.study_name in ["CIDR: The Genetic Basis of Aggressive Prostate Cancer: The Role of Rare Variation"]
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
 coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
  ORDER By samp.sample_id LIMIT 100
'@

$ws.Range("B3").Value = $newTumorQuery.TrimEnd("`r","`n")

# The extra leading blank line grows the wrapped-text row height.
$ws.Rows.Item(3).RowHeight = 204.75

# Move the active selection from D3 to B3.
$ws.Range("B3").Select()
